{"js": "// Appends two new work-log entries (1/4/24 and 2/4/24) to the first\n// table in the document, matching the author's commit:\n// \"Finalised design of Node-red flows from data decoding to display and\n// trigger events in PLC. Resolved issues with writing to Modbus due to\n// wrong data type. Dashboard displays info as expected and buttons work\n// as expected. Design report updated accordingly.\"\n\n// Helper: wrap a snippet of WordprocessingML block content (one or more\n// <w:p> elements) in a minimal OOXML \"flat package\" so it can be fed to\n// Range.insertOoxml() \u2014 this lets us control run-splitting exactly\n// (several runs per paragraph, <w:lastRenderedPageBreak/>, etc.) instead\n// of letting the host engine merge same-formatted runs together.\nfunction wrapOoxml(bodyInnerXml) {\n  return (\n    '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\"><pkg:xmlData>' +\n    '<Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\">' +\n    '<Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/>' +\n    '</Relationships></pkg:xmlData></pkg:part>' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body>' +\n    bodyInnerXml +\n    '</w:body></w:document>' +\n    '</pkg:xmlData></pkg:part></pkg:package>'\n  );\n}\n\n// The two new rows' \"Activity\" cell content, reproduced run-for-run\n// from the target OOXML diff.\nconst activity1Xml =\n  '<w:p xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:r><w:t>Worked on Node-red to decode data in flows and then use that to write to the PLC via Modbus. Updated design report to include th</w:t></w:r>' +\n  '<w:r><w:t xml:space=\"preserve\">e </w:t></w:r>' +\n  '<w:r><w:lastRenderedPageBreak/><w:t>code and scaling used</w:t></w:r>' +\n  '<w:r><w:t xml:space=\"preserve\">. Some issues during testing involving writing to Modbus and requires further troubleshooting </w:t></w:r>' +\n  '</w:p>';\n\nconst time1Xml =\n  '<w:p xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:r><w:lastRenderedPageBreak/><w:t>8</w:t></w:r>' +\n  '</w:p>';\n\nconst activity2Xml =\n  '<w:p xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:r><w:t>Found the i</w:t></w:r>' +\n  '<w:r><w:t>ssue</w:t></w:r>' +\n  '<w:r><w:t xml:space=\"preserve\"> of</w:t></w:r>' +\n  '<w:r><w:t xml:space=\"preserve\"> writing to Modbus</w:t></w:r>' +\n  '<w:r><w:t xml:space=\"preserve\"> was</w:t></w:r>' +\n  '<w:r><w:t xml:space=\"preserve\"> due to the assigned data type </w:t></w:r>' +\n  '<w:r><w:t>of the payload</w:t></w:r>' +\n  '<w:r><w:t xml:space=\"preserve\"> \\u2013 needs to be number 0 or 1</w:t></w:r>' +\n  '<w:r><w:t>. Finalisation of Node-red flow to UI dashboard.</w:t></w:r>' +\n  '</w:p>';\n\n// The diff does not specify a time value for the 2/4/24 row, so leave\n// its Time cell as a clean, empty paragraph (no placeholder run).\nconst time2Xml =\n  '<w:p xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"/>';\n\n// Locate the first table (Date / Activity / Time log).\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst logTable = tables.items[0];\n\n// Append two placeholder rows at the end of the table \u2014 plain text first,\n// then swap the Activity (and, for row 1, Time) cell contents for the\n// exact multi-run OOXML above.\nlogTable.addRows(\"End\", 2, [\n  [\"1/4/24\", \"\", \"8\"],\n  [\"2/4/24\", \"\", \"\"]\n]);\nawait context.sync();\n\nconst rows = logTable.rows;\nrows.load(\"items\");\nawait context.sync();\n\nconst newRow1 = rows.items[rows.items.length - 2];\nconst newRow2 = rows.items[rows.items.length - 1];\n\nconst row1Cells = newRow1.cells;\nrow1Cells.load(\"items\");\nconst row2Cells = newRow2.cells;\nrow2Cells.load(\"items\");\nawait context.sync();\n\nconst row1Activity = row1Cells.items[1].body;\nconst row1Time = row1Cells.items[2].body;\nconst row2Activity = row2Cells.items[1].body;\nconst row2Time = row2Cells.items[2].body;\n\nrow1Activity.insertOoxml(wrapOoxml(activity1Xml), \"Replace\");\nrow1Time.insertOoxml(wrapOoxml(time1Xml), \"Replace\");\nrow2Activity.insertOoxml(wrapOoxml(activity2Xml), \"Replace\");\nrow2Time.insertOoxml(wrapOoxml(time2Xml), \"Replace\");\nawait context.sync();\n", "ps1": "# Appends two new work-log entries (1/4/24 and 2/4/24) to the first\n# table in the document, matching the author's commit:\n# \"Finalised design of Node-red flows from data decoding to display and\n# trigger events in PLC. Resolved issues with writing to Modbus due to\n# wrong data type. Dashboard displays info as expected and buttons work\n# as expected. Design report updated accordingly.\"\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# Add two fresh rows at the end of the log table.\n$t.Rows.Add() | Out-Null\n$t.Rows.Add() | Out-Null\n\n$rowCount = $t.Rows.Count\n$row1Index = $rowCount - 1\n$row2Index = $rowCount\n\n$wRef = \"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"\n\n# --- Row 1: 1/4/24 --------------------------------------------------\n$t.Cell($row1Index, 1).Range.Text = \"1/4/24\"\n\n$activity1Xml = \"<w:p xmlns:w='$wRef'>\" +\n  \"<w:r><w:t>Worked on Node-red to decode data in flows and then use that to write to the PLC via Modbus. Updated design report to include th</w:t></w:r>\" +\n  \"<w:r><w:t xml:space='preserve'>e </w:t></w:r>\" +\n  \"<w:r><w:lastRenderedPageBreak/><w:t>code and scaling used</w:t></w:r>\" +\n  \"<w:r><w:t xml:space='preserve'>. Some issues during testing involving writing to Modbus and requires further troubleshooting </w:t></w:r>\" +\n  \"</w:p>\"\n$t.Cell($row1Index, 2).Range.InsertXML($activity1Xml) | Out-Null\n\n$time1Xml = \"<w:p xmlns:w='$wRef'><w:r><w:lastRenderedPageBreak/><w:t>8</w:t></w:r></w:p>\"\n$t.Cell($row1Index, 3).Range.InsertXML($time1Xml) | Out-Null\n\n# --- Row 2: 2/4/24 --------------------------------------------------\n$t.Cell($row2Index, 1).Range.Text = \"2/4/24\"\n\n$activity2Xml = \"<w:p xmlns:w='$wRef'>\" +\n  \"<w:r><w:t>Found the i</w:t></w:r>\" +\n  \"<w:r><w:t>ssue</w:t></w:r>\" +\n  \"<w:r><w:t xml:space='preserve'> of</w:t></w:r>\" +\n  \"<w:r><w:t xml:space='preserve'> writing to Modbus</w:t></w:r>\" +\n  \"<w:r><w:t xml:space='preserve'> was</w:t></w:r>\" +\n  \"<w:r><w:t xml:space='preserve'> due to the assigned data type </w:t></w:r>\" +\n  \"<w:r><w:t>of the payload</w:t></w:r>\" +\n  \"<w:r><w:t xml:space='preserve'> &#8211; needs to be number 0 or 1</w:t></w:r>\" +\n  \"<w:r><w:t>. Finalisation of Node-red flow to UI dashboard.</w:t></w:r>\" +\n  \"</w:p>\"\n$t.Cell($row2Index, 2).Range.InsertXML($activity2Xml) | Out-Null\n\n# The diff does not specify a time value for the 2/4/24 row, so leave\n# its Time cell as a clean, empty paragraph (no placeholder run).\n$time2Xml = \"<w:p xmlns:w='$wRef'/>\"\n$t.Cell($row2Index, 3).Range.InsertXML($time2Xml) | Out-Null\n"}
